# The author edited row 6 of the "bets" sheet: the 5th bet changed from a
# win to a loss. Column C holds the win/loss flag (1 = win, 0 = loss) and
# column E holds the gain/loss amount for that bet. Columns F and M are
# formula-driven (running balance and % change vs. the starting pot) and
# will recalculate automatically once the inputs change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

$ws.Range("C6").Value = 0
$ws.Range("E6").Value = -15250
